$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / Link / Price / Volume(1h) cell updates from the latest crypto data refresh.
# Price values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the source inlineStr cells) instead of auto-coercing to numbers and
# dropping formatting like trailing zeros (e.g. "27.00" -> 27).
$ws.Range("D2").Value = "'67.076.54"
$ws.Range("E2").Value = "  -3.62%  "
$ws.Range("D3").Value = "'3.517.86"
$ws.Range("E3").Value = "  -4.34%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'609.76"
$ws.Range("E5").Value = "  -5.82%  "
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("D7").Value = "'3.515.89"
$ws.Range("E7").Value = "  -4.30%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("E13").Value = "  -5.49%  "
$ws.Range("D14").Value = "'4.114.42"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "'31.61"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").Value = "'3.519.29"
$ws.Range("E16").Value = "  -4.20%  "
$ws.Range("D17").Value = "'67.008.84"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'15.32"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("D21").Value = "'443.89"
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("D22").Value = "'9.14"
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("D24").Value = "'77.74"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'3.660.04"
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  -7.08%  "
$ws.Range("D29").Value = "'8.15"
$ws.Range("E29").Value = "  -10.42%  "
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").Value = "'1.65"
$ws.Range("E31").Value = "  -4.23%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.159"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'25.68"
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("E35").Value = "  -4.83%  "
$ws.Range("E36").Value = "  -7.31%  "
$ws.Range("D37").Value = "'3.511.02"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("D38").Value = "'8.04"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'173.18"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "'2.15"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("D44").Value = "'0.0858"
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("D46").Value = "'45.23"
$ws.Range("E46").Value = "  -4.15%  "
$ws.Range("D47").Value = "'27.00"
$ws.Range("E47").Value = "  -6.71%  "
$ws.Range("E48").Value = "  -5.77%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("E51").Value = "  -5.24%  "
